$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF3").Value = "Dichotic_and_AFACT"
$ws.Range("AF5").Value = "MAB_and_AFACT"
$ws.Range("AF9").Value = "dichotic_phase"
$ws.Range("AF10").Value = "MAB_and_Digit_after"
$ws.Range("AF12").Value = "Dichotic_and_AFACT"
$ws.Range("AF13").Value = "MAB_phase"
$ws.Range("AF14").Value = "dichotic_phase"
$ws.Range("AF15").Value = "dichotic_phase"
$ws.Range("AF16").Value = "Dichotic_and_AFACT"
$ws.Range("AF17").Value = "Dichotic_and_AFACT"
$ws.Range("AF18").Value = "dichotic_phase"
$ws.Range("AF19").Value = "Digit_before_and_AFACT"
$ws.Range("AF20").Value = "Dichotic_and_AFACT"
$ws.Range("AF21").Value = "dichotic_phase"
$ws.Range("AF22").Value = "Digit_before_and_AFACT"
$ws.Range("AF23").Value = "MAB_and_Digit_after"
$ws.Range("AF24").Value = "MAB_phase"
$ws.Range("AF26").Value = "Digit_before_and_AFACT"
$ws.Range("AF27").Value = "MAB_and_AFACT"
$ws.Range("AF28").Value = "MAB_and_AFACT"
$ws.Range("AF29").Value = "MAB_and_Digit_after"
$ws.Range("AF30").Value = "MAB_phase"
$ws.Range("AF32").Value = "MAB_phase"
$ws.Range("AF34").Value = "dichotic_phase"
$ws.Range("AF35").Value = "Dichotic_and_AFACT"
$ws.Range("AF37").Value = "dichotic_phase"
$ws.Range("AF38").Value = "Dichotic_and_AFACT"
$ws.Range("AF39").Value = "MAB_and_AFACT"
$ws.Range("AF40").Value = "MAB_phase"
$ws.Range("AF41").Value = "Digit_before_and_AFACT"
$ws.Range("AF43").Value = "Dichotic_and_AFACT"
$ws.Range("AF44").Value = "dichotic_phase"
$ws.Range("AF45").Value = "dichotic_phase"
$ws.Range("AF46").Value = "MAB_phase"
$ws.Range("AF47").Value = "Dichotic_and_AFACT"
$ws.Range("AF48").Value = "Dichotic_and_AFACT"
$ws.Range("AF49").Value = "MAB_and_Digit_after"
$ws.Range("AF53").Value = "Digit_before_and_AFACT"
$ws.Range("AF54").Value = "MAB_phase"
$ws.Range("AF55").Value = "MAB_and_Digit_after"
$ws.Range("AF56").Value = "Dichotic_and_AFACT"
$ws.Range("AF59").Value = "Digit_before_and_AFACT"
$ws.Range("AF60").Value = "Dichotic_and_AFACT"
$ws.Range("AF61").Value = "dichotic_phase"
$ws.Range("AF62").Value = "dichotic_phase"
$ws.Range("AF64").Value = "dichotic_phase"
$ws.Range("AF65").Value = "Dichotic_and_AFACT"
$ws.Range("AF66").Value = "Digit_before_and_AFACT"
$ws.Range("AF67").Value = "MAB_and_AFACT"
$ws.Range("AF68").Value = "Dichotic_and_AFACT"
$ws.Range("AF69").Value = "MAB_phase"
$ws.Range("AF70").Value = "dichotic_phase"
$ws.Range("AF71").Value = "dichotic_phase"
$ws.Range("AF73").Value = "Digit_before_and_AFACT"
$ws.Range("AF74").Value = "MAB_phase"
$ws.Range("AF75").Value = "MAB_and_AFACT"
$ws.Range("AF76").Value = "MAB_and_Digit_after"
$ws.Range("AF77").Value = "Digit_before_and_AFACT"
$ws.Range("AF78").Value = "MAB_and_AFACT"
$ws.Range("AF79").Value = "MAB_phase"
$ws.Range("AF80").Value = "MAB_and_Digit_after"
$ws.Range("AF81").Value = "MAB_and_AFACT"